$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 (B:G) down into rows 3-11 (process bottom-up
# so we don't overwrite values before they've been copied).
for ($r = 10; $r -ge 2; $r--) {
    for ($col = 2; $col -le 7; $col++) {
        $srcCell = $ws.Cells.Item($r, $col)
        $dstCell = $ws.Cells.Item($r + 1, $col)
        $dstCell.Value = $srcCell.Value2
    }
}

# Write the new row of data into row 2 (B2:G2).
$ws.Cells.Item(2, 2).Value = -0.03550443442769693
$ws.Cells.Item(2, 3).Value = 0.4494482028570796
$ws.Cells.Item(2, 4).Value = 0.266414150275854
$ws.Cells.Item(2, 5).Value = 0.5161532236418309
$ws.Cells.Item(2, 6).Value = 0.5330038716810166
$ws.Cells.Item(2, 7).Value = 15
